$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New OHLCV rows appended after the existing last row (1084).
# Each entry: (row, Open-time, Open, High, Low, Close, Volume)
$newRows = @(
    @(1085, 45534.5, 4.272, 4.314, 4.101, 4.12, 904915.2),
    @(1086, 45534.66666666666, 4.121, 4.276, 4.103, 4.223, 653071.03),
    @(1087, 45534.83333333334, 4.222, 4.291, 4.222, 4.277, 206067.34),
    @(1088, 45535, 4.277, 4.313, 4.276, 4.293, 205668.45),
    @(1089, 45535.16666666666, 4.294, 4.318, 4.288, 4.289, 223582),
    @(1090, 45535.33333333334, 4.288, 4.299, 4.259, 4.287, 165031.54),
    @(1091, 45535.5, 4.286, 4.288, 4.25, 4.269, 192231.09),
    @(1092, 45535.66666666666, 4.27, 4.274, 4.218, 4.241, 206844.42),
    @(1093, 45535.83333333334, 4.24, 4.246, 4.228, 4.245, 53643.8),
    @(1094, 45536, 4.259, 4.259, 4.186, 4.219, 226796.4),
    @(1095, 45536.16666666666, 4.22, 4.225, 4.144, 4.203, 514495.51),
    @(1096, 45536.33333333334, 4.203, 4.206, 4.155, 4.192, 235809.43),
    @(1097, 45536.5, 4.193, 4.199, 4.108, 4.189, 551604.78),
    @(1098, 45536.66666666666, 4.186, 4.201, 4.142, 4.155, 111280.6),
    @(1099, 45536.83333333334, 4.195, 4.217, 4.032, 4.068, 529925.61),
    @(1100, 45537, 4.068, 4.126, 4.067, 4.115, 337326.93),
    @(1101, 45537.16666666666, 4.114, 4.133, 4.051, 4.073, 285375.06),
    @(1102, 45537.33333333334, 4.072, 4.209, 4.072, 4.19, 470763.11),
    @(1103, 45537.5, 4.19, 4.203, 4.177, 4.18, 48470.04),
    @(1104, 45537.66666666666, 4.157, 4.174, 4.135, 4.166, 245789.61),
    @(1105, 45537.83333333334, 4.167, 4.234, 4.163, 4.208, 338397.99),
    @(1106, 45538, 4.209, 4.279, 4.2, 4.259, 271316.83),
    @(1107, 45538.16666666666, 4.259, 4.268, 4.232, 4.265, 226666.16),
    @(1108, 45538.33333333334, 4.264, 4.275, 4.225, 4.265, 291531.36),
    @(1109, 45538.5, 4.264, 4.285, 4.096, 4.108, 656828.17),
    @(1110, 45538.66666666666, 4.107, 4.155, 4.096, 4.106, 340668.21),
    @(1111, 45538.83333333334, 4.107, 4.138, 4.056, 4.062, 258208.62),
    @(1112, 45539, 4.063, 4.097, 3.868, 4.038, 1457628.99),
    @(1113, 45539.16666666666, 4.037, 4.09, 4.032, 4.084, 396648),
    @(1114, 45539.33333333334, 4.085, 4.09, 4.012, 4.015, 253305.29),
    @(1115, 45539.5, 4.015, 4.162, 4.004, 4.157, 1208452.02),
    @(1116, 45539.66666666666, 4.156, 4.197, 4.095, 4.109, 486254.04),
    @(1117, 45539.83333333334, 4.123, 4.126, 4.122, 4.124, 468.33)
)

$lastRow = 1084
foreach ($r in $newRows) {
    $rowNum = $r[0]
    # Copy formatting (incl. the date style on column A) from the last existing row
    $ws.Range("A${lastRow}:F${lastRow}").Copy($ws.Range("A${rowNum}:F${rowNum}"))
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
}
